# Auto-generated script applying cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "  +0.18%  "
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "26.937.94"
$c.Style = "Normal"

$ws.Range("E3").Value = "  +1.26%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.883.73"
$c.Style = "Normal"

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("E5").Value = "  +0.10%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "305.07"
$c.Style = "Normal"

$ws.Range("E6").Value = "  -0.07%  "

$ws.Range("E7").Value = "  +1.86%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.5135"
$c.Style = "Normal"

$ws.Range("E8").Value = "  +2.43%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3733"
$c.Style = "Normal"

$ws.Range("E9").Value = "  +0.29%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.07190"
$c.Style = "Normal"

$ws.Range("E10").Value = "  +1.70%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "21.03"
$c.Style = "Normal"

$ws.Range("E11").Value = "  +0.68%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.8978"
$c.Style = "Normal"

$ws.Range("E12").Value = "  +2.20%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.07646"
$c.Style = "Normal"

$ws.Range("E13").Value = "  +0.29%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.869.36"
$c.Style = "Normal"

$ws.Range("E14").Value = "  -0.57%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "93.30"
$c.Style = "Normal"

$ws.Range("E15").Value = "  -0.06%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "5.224"
$c.Style = "Normal"

$ws.Range("E16").Value = "  -0.04%  "

$ws.Range("E17").Value = "  -0.34%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.000008465"
$c.Style = "Normal"

$ws.Range("E19").Value = "  -0.11%  "

$ws.Range("E20").Value = "  +0.24%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "26.997.58"
$c.Style = "Normal"

$ws.Range("E21").Value = "  +0.26%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "5.036"
$c.Style = "Normal"

$ws.Range("E22").Value = "  +0.31%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "2.120.10"
$c.Style = "Normal"

$ws.Range("E24").Value = "  -0.71%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "6.377"
$c.Style = "Normal"

$ws.Range("B25").Value = "LidoDAOToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("E25").Value = "  +9.89%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.284"
$c.Style = "Normal"

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("E26").Value = "  -0.92%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "146.31"
$c.Style = "Normal"

$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("E27").Value = "  -3.66%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "1.727"
$c.Style = "Normal"

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("E28").Value = "  +0.76%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "18.01"
$c.Style = "Normal"

$ws.Range("E29").Value = "  +0.68%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "113.86"
$c.Style = "Normal"

$ws.Range("E30").Value = "  +4.81%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "4.891"
$c.Style = "Normal"

$ws.Range("E31").Value = "  +1.57%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "4.766"
$c.Style = "Normal"

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.09165"
$c.Style = "Normal"

$ws.Range("E33").Value = "  -1.71%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.05027"
$c.Style = "Normal"

$ws.Range("E34").Value = "  +6.84%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.230"
$c.Style = "Normal"

$ws.Range("E35").Value = "  +2.19%  "

$ws.Range("E36").Value = "  -0.55%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "2.963"
$c.Style = "Normal"

$ws.Range("E37").Value = "  -0.21%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "3.260"
$c.Style = "Normal"

$ws.Range("E38").Value = "  +0.75%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "2.576"
$c.Style = "Normal"

$ws.Range("E39").Value = "  +0.46%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.5572"
$c.Style = "Normal"

$ws.Range("E40").Value = "  -0.82%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.01982"
$c.Style = "Normal"

$ws.Range("E42").Value = "  +5.87%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "9.007"
$c.Style = "Normal"

$ws.Range("E43").Value = "  +1.08%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "6.600"
$c.Style = "Normal"

$ws.Range("E44").Value = "  +0.69%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "118.42"
$c.Style = "Normal"

$ws.Range("E45").Value = "  +1.91%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.1496"
$c.Style = "Normal"

$ws.Range("E46").Value = "  +2.41%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.4790"
$c.Style = "Normal"

$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("E47").Value = "  -0.10%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.Style = "Normal"

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("E48").Value = "  +0.83%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "10.10"
$c.Style = "Normal"

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.590"
$c.Style = "Normal"

$ws.Range("E50").Value = "  +2.36%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "37.53"
$c.Style = "Normal"

$ws.Range("E51").Value = "  +1.21%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "63.73"
$c.Style = "Normal"
